$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header F1 and copy formatting from E1 (bold, border, centered)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Fill time_taken values for data rows 2-86
$ws.Range("F2").Value = "2021-10-05 10:50:07.812053"
$ws.Range("F3").Value = "2021-10-05 10:50:07.812064"
$ws.Range("F4").Value = "2021-10-05 10:50:07.812067"
$ws.Range("F5").Value = "2021-10-05 10:50:07.812070"
$ws.Range("F6").Value = "2021-10-05 10:50:07.812073"
$ws.Range("F7").Value = "2021-10-05 10:50:07.812076"
$ws.Range("F8").Value = "2021-10-05 10:50:07.812079"
$ws.Range("F9").Value = "2021-10-05 10:50:07.812081"
$ws.Range("F10").Value = "2021-10-05 10:50:07.812084"
$ws.Range("F11").Value = "2021-10-05 10:50:07.812086"
$ws.Range("F12").Value = "2021-10-05 10:50:07.812089"
$ws.Range("F13").Value = "2021-10-05 10:50:07.812091"
$ws.Range("F14").Value = "2021-10-05 10:50:07.812094"
$ws.Range("F15").Value = "2021-10-05 10:50:07.812096"
$ws.Range("F16").Value = "2021-10-05 10:50:07.812099"
$ws.Range("F17").Value = "2021-10-05 10:50:07.812101"
$ws.Range("F18").Value = "2021-10-05 10:50:07.812104"
$ws.Range("F19").Value = "2021-10-05 10:50:07.812107"
$ws.Range("F20").Value = "2021-10-05 10:50:07.812109"
$ws.Range("F21").Value = "2021-10-05 10:50:07.812112"
$ws.Range("F22").Value = "2021-10-05 10:50:07.812114"
$ws.Range("F23").Value = "2021-10-05 10:50:07.812116"
$ws.Range("F24").Value = "2021-10-05 10:50:07.812119"
$ws.Range("F25").Value = "2021-10-05 10:50:07.812121"
$ws.Range("F26").Value = "2021-10-05 10:50:07.812124"
$ws.Range("F27").Value = "2021-10-05 10:50:07.812127"
$ws.Range("F28").Value = "2021-10-05 10:50:07.812129"
$ws.Range("F29").Value = "2021-10-05 10:50:07.812132"
$ws.Range("F30").Value = "2021-10-05 10:50:07.812134"
$ws.Range("F31").Value = "2021-10-05 10:50:07.812137"
$ws.Range("F32").Value = "2021-10-05 10:50:07.812139"
$ws.Range("F33").Value = "2021-10-05 10:50:07.812142"
$ws.Range("F34").Value = "2021-10-05 10:50:07.812145"
$ws.Range("F35").Value = "2021-10-05 10:50:07.812147"
$ws.Range("F36").Value = "2021-10-05 10:50:07.812150"
$ws.Range("F37").Value = "2021-10-05 10:50:07.812152"
$ws.Range("F38").Value = "2021-10-05 10:50:07.812155"
$ws.Range("F39").Value = "2021-10-05 10:50:07.812158"
$ws.Range("F40").Value = "2021-10-05 10:50:07.812160"
$ws.Range("F41").Value = "2021-10-05 10:50:07.812163"
$ws.Range("F42").Value = "2021-10-05 10:50:07.812166"
$ws.Range("F43").Value = "2021-10-05 10:50:07.812168"
$ws.Range("F44").Value = "2021-10-05 10:50:07.812171"
$ws.Range("F45").Value = "2021-10-05 10:50:07.812173"
$ws.Range("F46").Value = "2021-10-05 10:50:07.812176"
$ws.Range("F47").Value = "2021-10-05 10:50:07.812178"
$ws.Range("F48").Value = "2021-10-05 10:50:07.812181"
$ws.Range("F49").Value = "2021-10-05 10:50:07.812183"
$ws.Range("F50").Value = "2021-10-05 10:50:07.812186"
$ws.Range("F51").Value = "2021-10-05 10:50:07.812188"
$ws.Range("F52").Value = "2021-10-05 10:50:07.812191"
$ws.Range("F53").Value = "2021-10-05 10:50:07.812193"
$ws.Range("F54").Value = "2021-10-05 10:50:07.812196"
$ws.Range("F55").Value = "2021-10-05 10:50:07.812199"
$ws.Range("F56").Value = "2021-10-05 10:50:07.812201"
$ws.Range("F57").Value = "2021-10-05 10:50:07.812204"
$ws.Range("F58").Value = "2021-10-05 10:50:07.812206"
$ws.Range("F59").Value = "2021-10-05 10:50:07.812209"
$ws.Range("F60").Value = "2021-10-05 10:50:07.812211"
$ws.Range("F61").Value = "2021-10-05 10:50:07.812213"
$ws.Range("F62").Value = "2021-10-05 10:50:07.812216"
$ws.Range("F63").Value = "2021-10-05 10:50:07.812218"
$ws.Range("F64").Value = "2021-10-05 10:50:07.812221"
$ws.Range("F65").Value = "2021-10-05 10:50:07.812223"
$ws.Range("F66").Value = "2021-10-05 10:50:07.812227"
$ws.Range("F67").Value = "2021-10-05 10:50:07.812230"
$ws.Range("F68").Value = "2021-10-05 10:50:07.812233"
$ws.Range("F69").Value = "2021-10-05 10:50:07.812235"
$ws.Range("F70").Value = "2021-10-05 10:50:07.812238"
$ws.Range("F71").Value = "2021-10-05 10:50:07.812240"
$ws.Range("F72").Value = "2021-10-05 10:50:07.812243"
$ws.Range("F73").Value = "2021-10-05 10:50:07.812245"
$ws.Range("F74").Value = "2021-10-05 10:50:07.812248"
$ws.Range("F75").Value = "2021-10-05 10:50:07.812250"
$ws.Range("F76").Value = "2021-10-05 10:50:07.812253"
$ws.Range("F77").Value = "2021-10-05 10:50:07.812255"
$ws.Range("F78").Value = "2021-10-05 10:50:07.812259"
$ws.Range("F79").Value = "2021-10-05 10:50:07.812262"
$ws.Range("F80").Value = "2021-10-05 10:50:07.812265"
$ws.Range("F81").Value = "2021-10-05 10:50:07.812268"
$ws.Range("F82").Value = "2021-10-05 10:50:07.812270"
$ws.Range("F83").Value = "2021-10-05 10:50:07.812273"
$ws.Range("F84").Value = "2021-10-05 10:50:07.812275"
$ws.Range("F85").Value = "2021-10-05 10:50:07.812278"
$ws.Range("F86").Value = "2021-10-05 10:50:07.812280"

Write-Host "Completed"
